# Insert a new row at position 557 (pushes existing rows 557:678 down to 558:679)
# and populate it with the new record's data, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(557).Insert()

$ws.Range("A557").Value = 10
$ws.Range("B557").Value = "Vega Modelo de Temuco"
$ws.Range("C557").Value = "La Araucanía"
$ws.Range("D557").Value = 45244
$ws.Range("E557").Value = 9
$ws.Range("F557").Value = 100112040
$ws.Range("G557").Value = "Cilantro"
$ws.Range("H557").Value = "Sin especificar"
$ws.Range("I557").Value = "Primera"
$ws.Range("J557").Value = 58
$ws.Range("K557").Value = 6000
$ws.Range("L557").Value = 6000
$ws.Range("M557").Value = 6000
$ws.Range("N557").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O557").Value = "Provincia de Cautín"
$ws.Range("P557").Value = 3000
$ws.Range("Q557").Value = 2
$ws.Range("R557").Value = "Hortaliza"
